# Casos VIHSIDA 2018 AL 2023.xlsx - "8 indicadores cap 2" edit
# Adds a yearly summary side-table (columns L:Q, rows 6-11) to the
# "VIHSIDA 2018 al 2022" sheet, and switches the active sheet/selection
# from sheet 1 to sheet 2.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("VIHSIDA 2018 al 2023")
$ws2 = $wb.Worksheets.Item("VIHSIDA 2018 al 2022")

# ---------------------------------------------------------------------
# New summary table on "VIHSIDA 2018 al 2022": headers in M6:Q6, then
# one row per year (2018-2022) in L7:Q11.
# ---------------------------------------------------------------------

# Header row (row 6) - reuse the existing header labels (shared strings)
$ws2.Range("M6").Value = "Mestizo Ladino"
$ws2.Range("N6").Value = "Maya"
$ws2.Range("O6").Value = "Xinca"
$ws2.Range("P6").Value = "Garífuna"
$ws2.Range("Q6").Value = "Otros "

# copy the formatting from the matching existing header cells (C6:G6)
$ws2.Range("C6:G6").Copy()
[void]$ws2.Range("M6:Q6").PasteSpecial(-4122)

# Year rows 2018-2022
$ws2.Range("L7").Value = 2018
$ws2.Range("M7").Value = 112
$ws2.Range("N7").Value = 68
$ws2.Range("O7").Value = 0
$ws2.Range("P7").Value = 0
$ws2.Range("Q7").Value = 0

$ws2.Range("L8").Value = 2019
$ws2.Range("M8").Value = 123
$ws2.Range("N8").Value = 51
$ws2.Range("O8").Value = 0
$ws2.Range("P8").Value = 0
$ws2.Range("Q8").Value = 0

$ws2.Range("L9").Value = 2020
$ws2.Range("M9").Value = 111
$ws2.Range("N9").Value = 35
$ws2.Range("O9").Value = 0
$ws2.Range("P9").Value = 0
$ws2.Range("Q9").Value = 0

$ws2.Range("L10").Value = 2021
$ws2.Range("M10").Value = 105
$ws2.Range("N10").Value = 30
$ws2.Range("O10").Value = 0
$ws2.Range("P10").Value = 1
$ws2.Range("Q10").Value = 0

$ws2.Range("L11").Value = 2022
$ws2.Range("M11").Value = 83
$ws2.Range("N11").Value = 20
$ws2.Range("O11").Value = 0
$ws2.Range("P11").Value = 0
$ws2.Range("Q11").Value = 0

# copy the formatting for the new "year" cells (L7:L11) from the existing
# year-label cell (B5), and the data cells (M7:Q11) from the matching
# existing data row (C7:G7)
$ws2.Range("B5").Copy()
[void]$ws2.Range("L7:L11").PasteSpecial(-4122)

$ws2.Range("C7:G7").Copy()
[void]$ws2.Range("M7:Q7").PasteSpecial(-4122)
[void]$ws2.Range("M8:Q8").PasteSpecial(-4122)
[void]$ws2.Range("M9:Q9").PasteSpecial(-4122)
[void]$ws2.Range("M10:Q10").PasteSpecial(-4122)
[void]$ws2.Range("M11:Q11").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Switch the active sheet/tab from "VIHSIDA 2018 al 2023" to
# "VIHSIDA 2018 al 2022", and update each sheet's remembered selection.
# ---------------------------------------------------------------------

[void]$ws1.Range("B2:H2").Select()

$ws2.Activate()
[void]$ws2.Range("L6:Q11").Select()
